$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# S28/G03: Portfolio backtest v1 (target weights) - mark as implemented
$ws.Range("D4").Value = "implemented"
$ws.Range("E4").Value = "Portfolio backtest v1 implemented: Target weights (EOD) with cadence/constraints/costs, equity+drawdown+actions."
$ws.Range("F4").Value = "27/12/2025 03:04"
